$d = $word.ActiveDocument

$lsq = [char]0x2018
$rsq = [char]0x2019
$ldq = [char]0x201C
$rdq = [char]0x201D

# Edit 1: shorten "...these nodes that it will take later. Note that this on its
# own can have issues...also empty." down to "...these nodes that it can take later."
$find1 = " that it will take later. Note that this on its own can have issues, deferring nodes means that they don" + $rsq + "t end up in the heap, and are therefore searched out of order. To prevent generating suboptimal paths, we must only perform this optimization if the node diagonally " + $ldq + "ahead" + $rdq + " of us is also empty."
$replace1 = " that it can take later."

$d.Content.Find.Execute($find1, $true, $false, $false, $false, $false, $true, 1, $false, $replace1, 2) | Out-Null

# Edit 2: "the node to our lower right will search it later" ->
# "the node to our lower left can search it later"
$find2 = "the node to our lower right will search it later"
$replace2 = "the node to our lower left can search it later"

$d.Content.Find.Execute($find2, $true, $false, $false, $false, $false, $true, 1, $false, $replace2, 2) | Out-Null

Write-Output "done"
